# Change the table style (Table Styles gallery selection) applied to the
# "Sources of finance" table on slide 6 from the default "Table_0" style
# ({26460AF4-DB96-4C0A-B356-4520C78A3495}) to the built-in style
# {42FD900B-3249-4BE7-A0D4-3FF618E997B0}.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(6)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $table = $shape.Table
        $table.ApplyStyle("{42FD900B-3249-4BE7-A0D4-3FF618E997B0}")
    }
}
